$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1. Clear the "new field" (Bad/red) highlight style from the tail of the
#        facebook row (row 5, cols AS:BG) *before* inserting the new row below
#        it, so the inserted row doesn't inherit that formatting.
$ws.Range("AS5:BG5").Style = "Normal"

# --- 2. Insert a new row above the old row 6 (facebook_posts) for the new
#        "facebook_videos" table header. This shifts every row from 6 down
#        (old 6 -> 7, ... old 17 -> 18) automatically.
$ws.Rows(6).Insert()

# --- 3. Populate the new row 6 with the facebook_videos header fields.
#        Values are written in the same order the strings were appended to
#        the shared-string table by the original author's edit: table name
#        first, then columns H..Q, then E..G, then D, then B/C last.
$ws.Cells.Item(6, 1).Value = "facebook_videos"

$ws.Cells.Item(6, 8).Value = "post_video_likes_by_reaction_type"
$ws.Cells.Item(6, 9).Value = "post_video_avg_time_watched"
$ws.Cells.Item(6, 10).Value = "post_video_social_actions"
$ws.Cells.Item(6, 11).Value = "post_video_view_time"
$ws.Cells.Item(6, 12).Value = "post_impressions_unique"
$ws.Cells.Item(6, 13).Value = "blue_reels_play_count"
$ws.Cells.Item(6, 14).Value = "fb_reels_total_plays"
$ws.Cells.Item(6, 15).Value = "fb_reels_replay_count"
$ws.Cells.Item(6, 16).Value = "post_video_retention_graph"
$ws.Cells.Item(6, 17).Value = "post_video_followers"

$ws.Cells.Item(6, 5).Value = "post_video_likes_by_reaction_type.REACTION_LIKE"
$ws.Cells.Item(6, 6).Value = "post_video_likes_by_reaction_type.REACTION_LOVE"
$ws.Cells.Item(6, 7).Value = "post_video_social_actions.COMMENT"

$ws.Cells.Item(6, 4).Value = "description"

$ws.Cells.Item(6, 2).Value = "created_time"
$ws.Cells.Item(6, 3).Value = "post_id"

# --- 4. Update the saved selection and print orientation.
$ws.Range("D21").Select()
$ws.PageSetup.Orientation = 1
